$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '70.842.12'
$ws.Range("E2").Value = '  +5.98%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.640.32'
$ws.Range("E3").Value = '  +5.91%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '593.61'
$ws.Range("E5").Value = '  +2.27%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '195.16'
$ws.Range("E6").Value = '  +3.51%  '
$ws.Range("E7").Value = '  +2.87%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.634.79'
$ws.Range("E8").Value = '  +5.90%  '
$ws.Range("E9").Value = '  -0.01%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.185'
$ws.Range("E10").Value = '  +8.19%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.681'
$ws.Range("E11").Value = '  +5.81%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '58.16'
$ws.Range("E12").Value = '  +1.73%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000308'
$ws.Range("E13").Value = '  +11.58%  '
$ws.Range("E14").Value = '  +5.80%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.222.99'
$ws.Range("E15").Value = '  +5.80%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.51'
$ws.Range("E16").Value = '  +9.03%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.638.77'
$ws.Range("E17").Value = '  +5.16%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '70.810.39'
$ws.Range("E18").Value = '  +6.00%  '
$ws.Range("E19").Value = '  +5.99%  '
$ws.Range("E20").Value = '  +2.52%  '
$ws.Range("E21").Value = '  +4.21%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '489.22'
$ws.Range("E22").Value = '  +2.73%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '19.36'
$ws.Range("E23").Value = '  +13.32%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.27'
$ws.Range("E24").Value = '  -2.38%  '
$ws.Range("E25").Value = '  +3.49%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '91.47'
$ws.Range("E26").Value = '  +2.61%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.19'
$ws.Range("E27").Value = '  +6.79%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.47'
$ws.Range("E28").Value = '  +5.27%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.60'
$ws.Range("E29").Value = '  +6.75%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.94'
$ws.Range("E30").Value = '  +6.84%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '32.85'
$ws.Range("E31").Value = '  +5.71%  '
$ws.Range("E32").Value = '  +10.16%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '12.31'
$ws.Range("E33").Value = '  +4.81%  '
$ws.Range("B34").Value = 'OKB'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '66.41'
$ws.Range("E34").Value = '  +2.84%  '
$ws.Range("B35").Value = 'Bittensor'
$ws.Range("C35").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '614.10'
$ws.Range("E35").Value = '  +2.52%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '40.53'
$ws.Range("E36").Value = '  +9.57%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0' + [string][char]0x2083 + '0833'
$ws.Range("E37").Value = '  +11.26%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.413'
$ws.Range("E38").Value = '  +6.11%  '
$ws.Range("E39").Value = '  +1.53%  '
$ws.Range("E40").Value = '  -0.13%  '
$ws.Range("E41").Value = '  +2.94%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.327.87'
$ws.Range("E42").Value = '  +4.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.19'
$ws.Range("E43").Value = '  +9.89%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.18'
$ws.Range("E44").Value = '  +16.88%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.83'
$ws.Range("E45").Value = '  +10.11%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0460'
$ws.Range("E46").Value = '  +7.07%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.71'
$ws.Range("E47").Value = '  +12.83%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.34'
$ws.Range("E48").Value = '  +3.07%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.140'
$ws.Range("E49").Value = '  +3.89%  '
$ws.Range("E50").Value = '  +2.08%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.999'
$ws.Range("E51").Value = '  -0.18%  '
